$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.874.61'
$ws.Range('E2').Value = '  -3.06%  '
$ws.Range('D3').Value = '2.918.27'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '584.47'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = '145.49'
$ws.Range('E6').Value = '  -4.75%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.505'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '2.917.43'
$ws.Range('E9').Value = '  -3.53%  '
$ws.Range('D10').Value = '6.86'
$ws.Range('E10').Value = '  +4.55%  '
$ws.Range('D11').Value = '0.145'
$ws.Range('E11').Value = '  -4.06%  '
$ws.Range('E12').Value = '  -3.84%  '
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').Value = '33.67'
$ws.Range('E14').Value = '  -5.20%  '
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '3.400.06'
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = '60.825.41'
$ws.Range('E17').Value = '  -3.09%  '
$ws.Range('D18').Value = '6.75'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').Value = '2.916.49'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('D20').Value = '431.37'
$ws.Range('E20').Value = '  -4.62%  '
$ws.Range('D21').Value = '13.63'
$ws.Range('E21').Value = '  -4.31%  '
$ws.Range('D22').Value = '0.683'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').Value = '  -4.50%  '
$ws.Range('D24').Value = '80.38'
$ws.Range('E24').Value = '  -3.28%  '
$ws.Range('D25').Value = '10.90'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').Value = '11.93'
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('E32').Value = '  -2.50%  '
$ws.Range('D33').Value = '26.56'
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('D34').Value = '0.107'
$ws.Range('E34').Value = '  -2.42%  '
$ws.Range('D35').Value = '0.0₃0874'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('D37').Value = '5.66'
$ws.Range('E37').Value = '  -4.31%  '
$ws.Range('D38').Value = '3.05'
$ws.Range('E38').Value = '  -3.82%  '
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('D40').Value = '49.72'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -4.59%  '
$ws.Range('D42').Value = '8.66'
$ws.Range('E42').Value = '  -4.46%  '
$ws.Range('D43').Value = '0.295'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('D44').Value = '41.04'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').Value = '378.20'
$ws.Range('E45').Value = '  -4.32%  '
$ws.Range('D46').Value = '0.0350'
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('D47').Value = '2.681.44'
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('D48').Value = '132.32'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '24.57'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('E51').Value = '  -1.76%  '
